$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2060
$ws1.Range("F4").Value = 856
$ws1.Range("F5").Value = 1180
$ws1.Range("F6").Value = 353

# Sheet "全部类型" (4th sheet): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2060
$ws4.Range("F6").Value = 856
$ws4.Range("F7").Value = 1180
$ws4.Range("F8").Value = 353
